$d = $word.ActiveDocument

# Paragraph 5 is the existing "User story 1" description paragraph
# ("For the user, I would like to add a pause button ... stop moving.").
# Insert two new empty paragraphs right after it; both inherit paragraph 5's
# (non-bold) formatting since we keep reusing its Range as the anchor.
$anchor = $d.Paragraphs.Item(5).Range
$anchor.InsertParagraphAfter()
$anchor.InsertParagraphAfter()

# New paragraph 6: "User story 2:" heading, bold.
$heading = $d.Paragraphs.Item(6)
$heading.Range.Text = "User story 2:"
$heading.Range.Bold = 1

# New paragraph 7: the user story body text, not bold.
$body = $d.Paragraphs.Item(7)
$body.Range.Text = "For the user, I would like to add a save button to save the layout and customization of a city. Priority is High and estimate time is 2 days. For testing, press the save button and save the file of the city with name."

# Append a new trailing empty paragraph after the bookmark paragraph,
# before the section break. Do this before touching the bookmark
# paragraph's own formatting so the new paragraph stays plain/unformatted.
$d.Paragraphs.Add() | Out-Null

# The trailing paragraph that holds the _GoBack bookmark now needs an
# explicit en-US language tag on its paragraph mark (matches the other
# paragraphs in the document).
$bookmarkPara = $d.Paragraphs.Item(8)
$bookmarkPara.Range.LanguageID = "en-US"
